# Regenerate merged AHB files
# Two-part edit applied to rows 38-173 of the active sheet:
#  1) Every "group header" row (first row of a new field-name group) gets the
#     same formatting already used by the existing header rows further up the
#     sheet (e.g. row 2): columns A,C:V -> style with grey fill, column B ->
#     bold + grey fill.
#  2) Column L ("AENDERUNG" marker) is cleared (value + style) for every
#     processed data row, matching the grey/centered "empty" look already
#     used on row 2's L cell.
#  A handful of rows (53, 95, 106) already carry a special highlight style
#  and are left completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose A:V formatting must be promoted to the "group header" look.
$fullRows = @(38,42,48,52,59,65,67,70,74,78,82,85,89,100,104,108,117,123,145,168,171)

# Rows that must be left completely untouched (already specially highlighted).
$skipRows = @(53,95,106)

# Template ranges already carrying the target formatting.
$templateRow = $ws.Range("A2:V2")
$templateL = $ws.Range("L2")

foreach ($r in $fullRows) {
    $templateRow.Copy()
    $ws.Range("A$r`:V$r").PasteSpecial(-4122)
}

for ($r = 38; $r -le 173; $r++) {
    if ($skipRows -contains $r) {
        continue
    }
    $templateL.Copy()
    $ws.Range("L$r").PasteSpecial(-4122)
    $ws.Range("L$r").ClearContents()
}

$excel.CutCopyMode = 0

Write-Host "done"
